# Generate Report for Archive
# - Status text "Ready for handoff" -> "In Translation" on every sheet that
#   shows it (Overview!E2:F3 plus the per-language "Status" column on the
#   zh-cn / de-de sheets).
# - The now-shorter status text lets those "Status" columns get narrower,
#   so their widths are tightened to match.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# --- Overview sheet: zh-cn / de-de summary status columns (E & F) ---
$overview.Range("E2:F3").Value = $newStatus

# --- zh-cn / de-de sheets: "Status" table column (C) ---
$zhcn.Range("C2:C3").Value = $newStatus
$dede.Range("C2:C3").Value = $newStatus

# --- Narrow the affected status columns to fit the shorter text ---
$newWidth = 12.5

$overview.Columns.Item(5).ColumnWidth = $newWidth   # column E
$overview.Columns.Item(6).ColumnWidth = $newWidth   # column F
$zhcn.Columns.Item(3).ColumnWidth = $newWidth        # column C
$dede.Columns.Item(3).ColumnWidth = $newWidth        # column C
